$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number-looking string must keep their original
# "text" storage (matches the source inlineStr cells), so force text format first.
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

# Apply updated Price / Volume(1h) values, and row-content swaps for rows 44-47
$ws.Range('D2').Value = '57.974.07'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').Value = '2.447.82'
$ws.Range('E3').Value = '  +0.57%  '
$ws.Range('D4').Value = '0.995'
$ws.Range('E4').Value = '  -0.57%  '
$ws.Range('D5').Value = '507.78'
$ws.Range('E5').Value = '  -2.20%  '
$ws.Range('D6').Value = '133.51'
$ws.Range('E6').Value = '  +4.47%  '
$ws.Range('D7').Value = '0.995'
$ws.Range('E7').Value = '  -0.61%  '
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('D9').Value = '2.444.66'
$ws.Range('E9').Value = '  +0.15%  '
$ws.Range('E10').Value = '  +1.14%  '
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('E12').Value = '  +1.30%  '
$ws.Range('E13').Value = '  -5.83%  '
$ws.Range('D14').Value = '2.873.75'
$ws.Range('E14').Value = '  +0.16%  '
$ws.Range('D15').Value = '57.817.77'
$ws.Range('E15').Value = '  +0.48%  '
$ws.Range('D16').Value = '21.89'
$ws.Range('E16').Value = '  +2.20%  '
$ws.Range('D17').Value = '0.0000136'
$ws.Range('E17').Value = '  +3.41%  '
$ws.Range('D18').Value = '2.391.21'
$ws.Range('E18').Value = '  -2.16%  '
$ws.Range('D19').Value = '10.32'
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('D20').Value = '315.42'
$ws.Range('E20').Value = '  +1.74%  '
$ws.Range('E21').Value = '  +0.73%  '
$ws.Range('E22').Value = '  +6.22%  '
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('E24').Value = '  -1.44%  '
$ws.Range('E25').Value = '  +1.21%  '
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  -0.71%  '
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('E28').Value = '  -4.39%  '
$ws.Range('D29').Value = '7.56'
$ws.Range('E29').Value = '  +5.53%  '
$ws.Range('D30').Value = '170.76'
$ws.Range('E30').Value = '  -1.40%  '
$ws.Range('D31').Value = '0.0₃0735'
$ws.Range('E31').Value = '  +0.86%  '
$ws.Range('E32').Value = '  +0.98%  '
$ws.Range('D33').Value = '6.13'
$ws.Range('E33').Value = '  -0.10%  '
$ws.Range('E34').Value = '  +1.53%  '
$ws.Range('D35').Value = '0.997'
$ws.Range('E35').Value = '  -0.16%  '
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  -0.24%  '
$ws.Range('D37').Value = '18.08'
$ws.Range('E37').Value = '  +1.91%  '
$ws.Range('E38').Value = '  +4.99%  '
$ws.Range('E39').Value = '  +4.20%  '
$ws.Range('D40').Value = '36.76'
$ws.Range('E40').Value = '  +1.63%  '
$ws.Range('E41').Value = '  +3.47%  '
$ws.Range('E42').Value = '  +2.79%  '
$ws.Range('D43').Value = '134.92'
$ws.Range('E43').Value = '  +13.47%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '5.03'
$ws.Range('E44').Value = '  +6.92%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').Value = '3.39'
$ws.Range('E45').Value = '  +1.14%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').Value = '255.80'
$ws.Range('E46').Value = '  +2.70%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = '0.572'
$ws.Range('E47').Value = '  -1.03%  '
$ws.Range('D48').Value = '0.0916'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('E49').Value = '  +0.71%  '
$ws.Range('E50').Value = '  +2.72%  '
$ws.Range('D51').Value = '17.20'
$ws.Range('E51').Value = '  +2.06%  '
